# Auto-generated script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.452.45"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -4.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.291.71"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -6.17%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "180.88"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -9.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "527.56"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.603"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.287.46"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -6.13%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.614"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.70"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -5.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.133"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -6.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000260"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.10"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -7.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.826.92"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -5.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.303.97"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.117"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -5.10%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.283.03"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.28%  "
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.56"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.12"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -6.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.958"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "374.31"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.80"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.95"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.08"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -6.31%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.10"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.69"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.55"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -6.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.42"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.94"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -6.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "647.60"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.73"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.31"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.105"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.00"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -6.57%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.393"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.55"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.47%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0700"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.126"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.868.57"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -6.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.48"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.70"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -10.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0399"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.65"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.86"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +9.42%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.03"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +6.32%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.55"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -6.39%  "
